$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "qlksJ932"
$ws.Range("B2").Value = 23081808
$ws.Range("C2").Value = "grxajwi65"
$ws.Range("D2").Value = "c32!&JrX"
$ws.Range("F2").Value = "TidmbvVd"
$ws.Range("G2").Value = "frOZ"
